$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (Förändrad / "Changed") date values for rows 2-11 were bumped
# from serial 45183 (2023-09-14) to serial 45184 (2023-09-15).
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
